# The "Total Statistics" sheet (Sheet1) pulls RIP/test-suite counts from
# external linked workbooks via formulas such as =[4]Sheet1!$G$1 (row 5,
# the "Running" folder). Those source workbooks were updated through the
# new delete-RIPs API, so the cached numbers that flow into this summary
# change. We push the refreshed source numbers into the four cells that
# anchor that external link (B5, C5, G5, H5); every other cell on the
# sheet (N1, L2, N2, L3, N3, P3, E5, I5, L5, L6, L7, ...) is a normal
# in-workbook formula and recalculates automatically from these.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 5    # was 4  -> [4]Sheet1!$G$1
$ws.Range("C5").Value = 4    # was 0  -> [4]Sheet1!$G$2
$ws.Range("G5").Value = 10   # was 0  -> [4]Sheet1!$G$5
$ws.Range("H5").Value = 54   # was 55 -> [4]Sheet1!$G$4
